$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for rows 2,3,5,6,7
# These changes apply identically to the "展览" and "全部类型" sheets.
$updates = @{
    "F2" = 6994
    "F3" = 54
    "F5" = 82
    "F6" = 1079
    "F7" = 171
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
